# Replace the textual "<year> <month range>" labels in column A (rows 3-38)
# with plain numeric year values, matching the target edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2022
$ws.Range("A4").Value = 2021
$ws.Range("A5").Value = 2020
$ws.Range("A6").Value = 2019
$ws.Range("A7").Value = 2018
$ws.Range("A8").Value = 2017
$ws.Range("A9").Value = 2016
$ws.Range("A10").Value = 2015
$ws.Range("A11").Value = 2014
$ws.Range("A12").Value = 2013
$ws.Range("A13").Value = 2012
$ws.Range("A14").Value = 2012
$ws.Range("A15").Value = 2011
$ws.Range("A16").Value = 2011
$ws.Range("A17").Value = 2011
$ws.Range("A18").Value = 2010
$ws.Range("A19").Value = 2009
$ws.Range("A20").Value = 2009
$ws.Range("A21").Value = 2008
$ws.Range("A22").Value = 2007
$ws.Range("A23").Value = 2006
$ws.Range("A24").Value = 2005
$ws.Range("A25").Value = 2004
$ws.Range("A26").Value = 2003
$ws.Range("A27").Value = 2003
$ws.Range("A28").Value = 2002
$ws.Range("A29").Value = 2001
$ws.Range("A30").Value = 2001
$ws.Range("A31").Value = 2001
$ws.Range("A32").Value = 2000
$ws.Range("A33").Value = 1999
$ws.Range("A34").Value = 1998
$ws.Range("A35").Value = 1997
$ws.Range("A36").Value = 1997
$ws.Range("A37").Value = 1996
$ws.Range("A38").Value = 1995

# Move the active selection from A42 to A39, as in the target workbook.
$ws.Range("A39").Select()
